$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

# --- Row 1 header: drop the stray repeated "production_rate"-style header
#     cells in C1:F1 and relabel A1/B1 to the generic parameter/value header ---
$ws.Range("A1").Value = "optimization_parameter"
$ws.Range("B1").Value = "value"
$ws.Range("C1:F1").ClearContents()

# --- Row 8: "Model" label renamed to "production_function" (value unchanged) ---
$ws.Range("A8").Value = "production_function"

# --- Insert a new row for the "L_curve" flag right after the production
#     function row, pushing estimate_params..Sheet down by one ---
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").Value = 0
$ws.Range("B9").NumberFormat = "0.00E+00"

# --- Remove the obsolete "Deletion" row entirely (now sitting at row 17
#     after the insert above) ---
$ws.Rows.Item(17).Delete()

# --- This sheet becomes the active / selected tab, with C1:F1 highlighted ---
$ws.Activate()
$ws.Range("C1:F1").Select()
